$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the last bullet: "...rows that actually had values..." ->
#    "...rows that had values..." and collapse the old multi-run / proofErr
#    structure down to a single run (matching the target OOXML), while
#    keeping "Next, several " as its own separate leading run.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)

$oldTail = "lines of code were written to keep only the rows that actually had values in specified columns."
$newTail = "lines of code were written to keep only the rows that had values in specified columns."
$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newTail, 2) | Out-Null

# Re-split "Next, several " back out into its own run (the text-fix above
# merges the whole paragraph into one run). Toggling a character attribute
# on/off forces the engine to keep that span as a distinct run.
$p7Start = $p7.Range.Start
$splitAt = $p7Start + "Next, several ".Length
$leadRange = $d.Range($p7Start, $splitAt)
$leadRange.Bold = 1
$leadRange.Bold = 0

# ---------------------------------------------------------------------------
# 2) Add the two new bullet points after it, inheriting the same
#    ListParagraph / numPr formatting automatically via InsertParagraphAfter.
# ---------------------------------------------------------------------------
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Perform groupby with mean aggregation to look at average point and price for each title."

$p8Start = $p8.Range.Start
$p8End = $p8.Range.End
$len1 = "Perform groupby with mean ".Length
$len2 = "aggregation".Length
$len3 = " to ".Length
$off1 = $p8Start + $len1
$off2 = $p8Start + $len1 + $len2
$off3 = $p8Start + $len1 + $len2 + $len3

# Split from the rightmost boundary first so earlier offsets stay valid.
$s3 = $d.Range($off3, $p8End)
$s3.Bold = 1
$s3.Bold = 0

$s2 = $d.Range($off2, $off3)
$s2.Bold = 1
$s2.Bold = 0

$s1 = $d.Range($off1, $off2)
$s1.Bold = 1
$s1.Bold = 0

$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Perform sort_values to determine the Top 5 & Bottom 5 popular titles."

$p9Start = $p9.Range.Start
$p9End = $p9.Range.End
$len4 = "Perform sort_values to determine the Top 5 & Bottom 5 popular ".Length
$off4 = $p9Start + $len4

$s4 = $d.Range($off4, $p9End)
$s4.Bold = 1
$s4.Bold = 0

Write-Output "Done"
